$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 241, pushing existing rows 241:366 down to 243:368.
$ws.Rows.Item(241).Insert()
$ws.Rows.Item(241).Insert()

# The row that used to be at 241 is now at 243 (identical data: same market/product/quality
# row, just a different date). Duplicate it into the two freshly-inserted rows 241/242, then
# overwrite the handful of cells (Fecha, Volumen, Precio promedio ponderado, Precio $/Kg) that
# differ for the newly reported week.
$ws.Range("A243:R244").Copy()
$ws.Range("A241").PasteSpecial(-4104)
$excel.CutCopyMode = 0

$ws.Range("D241").Value = 44452
$ws.Range("J241").Value = 2100
$ws.Range("M241").Value = 7405
$ws.Range("P241").Value = 1234

$ws.Range("D242").Value = 44452
$ws.Range("J242").Value = 650

$wb.Save()
